$d = $word.ActiveDocument

$replacements = @(
    @("852÷3=284, 0", "282÷7=40, 2"),
    @("243÷5=48, 3", "697÷8=87, 1"),
    @("107÷8=13, 3", "560÷3=186, 2"),
    @("990÷4=247, 2", "168÷9=18, 6"),
    @("392÷8=49, 0", "181÷8=22, 5"),
    @("775÷6=129, 1", "257÷6=42, 5"),
    @("992÷7=141, 5", "541÷8=67, 5"),
    @("862÷2=431, 0", "177÷7=25, 2"),
    @("497÷5=99, 2", "731÷9=81, 2"),
    @("755÷2=377, 1", "125÷3=41, 2"),
    @("146÷8=18, 2", "268÷2=134, 0"),
    @("361÷7=51, 4", "722÷9=80, 2"),
    @("691÷9=76, 7", "474÷9=52, 6"),
    @("221÷2=110, 1", "532÷4=133, 0"),
    @("648÷9=72, 0", "522÷4=130, 2"),
    @("976÷3=325, 1", "583÷5=116, 3"),
    @("234÷8=29, 2", "656÷2=328, 0"),
    @("364÷3=121, 1", "530÷4=132, 2"),
    @("403÷8=50, 3", "476÷6=79, 2"),
    @("196÷4=49, 0", "829÷7=118, 3"),
    @("141÷3=47, 0", "102÷7=14, 4"),
    @("723÷2=361, 1", "275÷5=55, 0"),
    @("536÷6=89, 2", "829÷4=207, 1"),
    @("665÷6=110, 5", "993÷3=331, 0"),
    @("367÷7=52, 3", "723÷5=144, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
